# schedule_aller.xlsx - "add colors TOPAZE and AMETHYSTE"
#
# 1. Recolor the header row fill (was a blue 4472C4) to a dark gold/brown
#    (996515), and set the header font to Arial.
# 2. Recolor the two alternating data-row fills:
#      BDD7EE (light blue)  -> FFE5B4 (peach)
#      E2EFDA (light green) -> FFD700 (gold)
# 3. Translate the header labels to French / capitalise them.
# 4. Narrow column F by one character.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- helper: VBA-style RGB() packing (R + G*256 + B*65536) -------------
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$oldHeader = RGBVal 0x44 0x72 0xC4   # 4472C4
$oldOdd    = RGBVal 0xBD 0xD7 0xEE   # BDD7EE
$oldEven   = RGBVal 0xE2 0xEF 0xDA   # E2EFDA

$newHeader = RGBVal 0x99 0x65 0x15   # 996515 "TOPAZE"
$newOdd    = RGBVal 0xFF 0xE5 0xB4   # FFE5B4
$newEven   = RGBVal 0xFF 0xD7 0x00   # FFD700 "AMETHYSTE"

# ---- 1. Header row (A1:G1): new fill colour + Arial font ---------------
$headerRange = $ws.Range("A1:G1")
$headerRange.Interior.Color = $newHeader
$headerRange.Interior.PatternColor = $newHeader
$headerRange.Font.Name = "Arial"

# ---- 2. Recolour every data row's fill based on its current colour -----
$dim = $ws.UsedRange
$lastRow = $dim.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $rowRange = $ws.Range("A" + $r + ":G" + $r)
    $current = $rowRange.Cells.Item(1, 1).Interior.Color

    if ($current -eq $oldOdd) {
        $rowRange.Interior.Color = $newOdd
        $rowRange.Interior.PatternColor = $newOdd
    } elseif ($current -eq $oldEven) {
        $rowRange.Interior.Color = $newEven
        $rowRange.Interior.PatternColor = $newEven
    }
}

# ---- 3. Header text translations ---------------------------------------
$ws.Range("A1").Value = "Round"
$ws.Range("B1").Value = "Début"
$ws.Range("C1").Value = "Fin"
$ws.Range("D1").Value = "Équipe 1"
$ws.Range("E1").Value = "Équipe 2"
$ws.Range("F1").Value = "Durée"
$ws.Range("G1").Value = "Phase"

# ---- 4. Narrow column F (6) from 10 to 9 --------------------------------
$ws.Columns.Item(6).ColumnWidth = 9
